# Applies the "Wrote FRs for MSS for UC 5, 6, and 10" commit changes
# to the Requirements worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update/replace existing requirement rows (UC2 -> UC2/UC6 merge, UC3 renumbering) ---
$ws.Range("B18").Value = 'UC2A, UC6A'
$ws.Range("C18").Value = 'The system must be able to edit a post/tool-request details as the user specifies'
$ws.Range("D18").Value = 'The post/tool-request details will be updated in the database'
$ws.Range("E18").Value = '-'
$ws.Range("H18").Value = 'Created 4/3/21' + [char]10 + 'Edited 5/3/21'

$ws.Range("B19").Value = 'UC2B, UC6B'
$ws.Range("C19").Value = 'The system must be able to delete a post/tool-request '
$ws.Range("D19").Value = 'The post/tool-request will be removed from the database'
$ws.Range("E19").Value = '-'
$ws.Range("H19").Value = 'Created 4/3/21' + [char]10 + 'Edited 5/3/21'

$ws.Range("B20").Value = 'UC3,'
$ws.Range("C20").Value = 'The system must be able to send notifications to a user.'
$ws.Range("D20").Value = 'The user will recive a notification'

$ws.Range("B21").Value = 'UC3'
$ws.Range("C21").Value = 'The system must be able to determine if a user is authorized to view, accept, and deny tool-requests to a certain post.'
$ws.Range("D21").Value = 'an unauthorized user will not be able to view, accept, or deny a tool-request'

$ws.Range("B22").Value = 'UC3'
$ws.Range("C22").Value = 'The system must allow the user (owner) to view tool-requests sent to their posts'
$ws.Range("D22").Value = 'a screen with all the requests will be displayed'

$ws.Range("B23").Value = 'UC3'
$ws.Range("C23").Value = 'The system must allow the user (owner) to accept a tool-request'
$ws.Range("D23").Value = 'the request will be accepted and the tool marked as unavailable'

$ws.Range("B24").Value = 'UC3'
$ws.Range("C24").Value = 'The system must be able to open a chat room between the owner and renter'
$ws.Range("D24").Value = 'a chat room will be open where the owner and renter can send and recive messages from each other'

# --- New requirement rows for UC5 / UC6 (tool-requests) ---
$ws.Range("B26").Value = 'UC5, UC6'
$ws.Range("C26").Value = 'The system must be able to determine if a user is authorized to send, edit, or remove a tool-request to a tool-post.'
$ws.Range("D26").Value = 'The user''s request to add, edit or delete a tool-request will be denied if (s)he are not authorized to'
$ws.Range("H26").Value = 'Created 5/3/21'
$ws.Rows.Item(26).RowHeight = 60

$ws.Range("B27").Value = 'UC5'
$ws.Range("C27").Value = 'The system must allow the user (renter) to send a tool-request to a post.'
$ws.Range("D27").Value = 'a tool-request screen will be displayed'
$ws.Range("H27").Value = 'Created 5/3/21'
$ws.Rows.Item(27).RowHeight = 30

$ws.Range("B28").Value = 'UC5'
$ws.Range("C28").Value = 'The system must be able to add a tool-request to a tool-post'
$ws.Range("D28").Value = 'the tool-request will be added to the post''s requests, where the owner can accept, or deny it.'
$ws.Range("H28").Value = 'Created 5/3/21'
$ws.Rows.Item(28).RowHeight = 45

$ws.Range("B29").Value = 'UC6'
$ws.Range("C29").Value = 'The system must allow the user (renter) to edit or remove his/her tool-request'
$ws.Range("D29").Value = 'a screen to edit the request will be displayed with a button to delete the tool-request'
$ws.Range("H29").Value = 'Created 5/3/21'
$ws.Rows.Item(29).RowHeight = 45

# --- New requirement rows for UC10 (reviews / ratings) ---
$ws.Range("B31").Value = 'UC10'
$ws.Range("C31").Value = 'The system must be able to determine if two users had a previous rental-transaction.'
$ws.Range("D31").Value = 'the system will confirm if two users had a previous transaction or not.'
$ws.Range("H31").Value = 'Created 5/3/21'
$ws.Rows.Item(31).RowHeight = 45

$ws.Range("B32").Value = 'UC10A'
$ws.Range("C32").Value = 'The system must allow the user to add/edit a review of another user if they had a previous transaction'
$ws.Range("D32").Value = 'A form to submit the review details will be displayed'
$ws.Range("H32").Value = 'Created 5/3/21'
$ws.Rows.Item(32).RowHeight = 30

$ws.Range("B33").Value = 'UC10'
$ws.Range("C33").Value = 'The system must be able to add/edit/delete a review to a user''s profile'
$ws.Range("D33").Value = 'a review will be added/edited/deleted and the total rating score will be recalculated'
$ws.Range("H33").Value = 'Created 5/3/21'
$ws.Rows.Item(33).RowHeight = 60

$ws.Range("B34").Value = 'UC10B'
$ws.Range("C34").Value = 'The system must be able to determine if a user is authorized to delete a review'
$ws.Range("D34").Value = 'the delete request will be denied if the user isn''t authorized'
$ws.Range("H34").Value = 'Created 5/3/21'
$ws.Rows.Item(34).RowHeight = 30

$ws.Range("B35").Value = 'UC10B'
$ws.Range("C35").Value = 'The system must allow the user to delete a review.'
$ws.Range("D35").Value = 'an option (button) to delete the review will be displayed'
$ws.Range("H35").Value = 'Created 5/3/21'
$ws.Rows.Item(35).RowHeight = 30

# --- Update the active selection/view to match the authored state ---
$ws.Range("D39").Select()
